$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")
$ws.Rows("122").Delete()
